# LMS prod Configuration Set Up
# Updates the "LMSPROD" worksheet with new school/classroom/section
# identifiers and numeric IDs (stored as text, matching the source data).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LMSPROD")

# --- Plain text values (School / Classroom / Section names) ---------------
$ws.Range("A2").Value = "FPK12School95878"
$ws.Range("B2").Value = "FPK12Classroom27824"
$ws.Range("C2").Value = "FPK12Section13722"

# fpk12student row label that was previously blank
$ws.Range("D5").Value = "fpk12student"

# --- Numeric-looking IDs that must remain stored as TEXT -------------------
# Directly assigning a numeric-looking string (e.g. "43388") gets
# auto-coerced to a number by the Value setter, so build the text via a
# formula that evaluates to a string, then paste-special just the value
# back into place - this preserves the text type and the cell's
# existing style.
$scratch = $ws.Range("Z100")

$scratch.Formula = '="43388"'
$scratch.Copy()
$ws.Range("E3").PasteSpecial(-4163)

$scratch.Formula = '="9801"'
$scratch.Copy()
$ws.Range("E4").PasteSpecial(-4163)

$scratch.Formula = '="38240"'
$scratch.Copy()
$ws.Range("E5").PasteSpecial(-4163)

$scratch.Clear()
